# Edit script implementing:
#  1. Slide 2 (sldId 257): add a Title placeholder shape with text "Descripcion"
#  2. Slide 8 (sldId 263): append a period to the last sentence of the body text
#  3. New Slide 9 (sldId 264) appended at the end: "Conclusiones resultados" with
#     the two extra conclusion bullets.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 2 — add the missing Title shape ("Descripcion")
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$titleSource = $p.Slides.Item(8).Shapes.Item(1)
$titleSource.Copy()
$newTitle = $slide2.Shapes.Paste()
$newTitle.Left = 1451579 / 12700
$newTitle.Top = 804519 / 12700
$newTitle.Width = 9603275 / 12700
$newTitle.Height = 1049235 / 12700
$newTitle.TextFrame.TextRange.Text = "Descripcion"
$newTitle.TextFrame.TextRange.LanguageID = "es-MX"

# ---------------------------------------------------------------------------
# 2) Slide 8 — fix the trailing sentence so it ends with a period
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange
$lastPara = $body8.Paragraphs(3, 1)
$lastPara.Text = $lastPara.Text + "."

# ---------------------------------------------------------------------------
# 3) Add new Slide 9 — "Conclusiones resultados"
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Add($p.Slides.Count + 1, 2)

$slide9.Shapes.Item(1).Name = "Título 1"
$slide9.Shapes.Item(2).Name = "Marcador de contenido 2"

$title9 = $slide9.Shapes.Item(1).TextFrame.TextRange
$title9.Text = "Conclusiones resultados"
$title9.LanguageID = "es-MX"

$body9 = $slide9.Shapes.Item(2).TextFrame.TextRange
$bullet1 = "- La Fabricación de alimentos con aceites No Vegetales es más del doble que los Vegetales, esto nos dice que la mayor ganancias proviene del OIL por la función objetivo; lo que concluimos con esto es que la demanda del OIL debe de irse igualanado a la del aceite Vegetal para tener un consumo de alimentos más sano."
$bullet2 = "- Las ganancias máximas durante los primeros 6 meses del año no son tan buenas."
$body9.Text = $bullet1 + "`r" + "`r" + $bullet2
$body9.LanguageID = "es-MX"
